{"js": "// The \"type\" attribute is removed from the Person(...) entity in the\n// logical-model iteration-1 section:\n//   Person(person_id, fname, lname, phone, mail, bday,type)\n//     -> Person(person_id, fname, lname, phone, mail, bday)\n//\n// (The second Person(...) line further down, which also carries\n// \"address_id\", is a different iteration's entity and is left untouched.)\n\nconst body = context.document.body;\n\n// Word keeps a single \"_GoBack\" bookmark that always tracks the location of\n// the most recent edit. Before this edit it sits in the \"Messages(...)\"\n// paragraph; after this edit Word will have silently relocated it to sit\n// right where the user's cursor left off (immediately after \"bday\", where\n// \",type\" used to be). Remove it from its old spot first.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Locate the exact entity definition to edit and strip \",type\" from it.\nconst matches = body.search(\"Person(person_id, fname, lname, phone, mail, bday,type)\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  matches.items[0].insertText(\n    \"Person(person_id, fname, lname, phone, mail, bday)\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// Re-seat \"_GoBack\" right after \"bday\" in that same paragraph (the spot the\n// cursor sits once \",type\" is deleted), splitting the text into two runs\n// with the bookmark in between \u2014 matching Word's own edit tracking.\nconst bdayMatches = body.search(\"bday\", { matchCase: true });\nbdayMatches.load(\"items\");\nawait context.sync();\n\nif (bdayMatches.items.length > 0) {\n  const afterBday = bdayMatches.items[0].getRange(\"End\");\n  afterBday.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The \"type\" attribute is removed from the Person(...) entity in the\n# logical-model iteration-1 section:\n#   Person(person_id, fname, lname, phone, mail, bday,type)\n#     -> Person(person_id, fname, lname, phone, mail, bday)\n#\n# (The second Person(...) line further down, which also carries\n# \"address_id\", is a different iteration's entity and is left untouched.)\n\n$d = $word.ActiveDocument\n\n# Word keeps a single \"_GoBack\" bookmark that always tracks the location of\n# the most recent edit. Before this edit it sits in the \"Messages(...)\"\n# paragraph; after this edit Word will have silently relocated it to sit\n# right where the user's cursor left off (immediately after \"bday\", where\n# \",type\" used to be). Remove it from its old spot first.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n# Locate the exact entity definition to edit and strip \",type\" from it.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Person(person_id, fname, lname, phone, mail, bday,type)\"\n$find.Execute()\nif ($find.Found) {\n  $rng.Text = \"Person(person_id, fname, lname, phone, mail, bday)\"\n}\n\n# Re-seat \"_GoBack\" right after \"bday\" in that same paragraph (the spot the\n# cursor sits once \",type\" is deleted), splitting the text into two runs\n# with the bookmark in between - matching Word's own edit tracking.\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.Text = \"bday\"\n$find2.Execute()\nif ($find2.Found) {\n  $afterBday = $d.Range($rng2.End, $rng2.End)\n  $d.Bookmarks.Add(\"_GoBack\", $afterBday)\n}\n"}
